$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "KPIs" - refreshed top-line numbers (revenue dropped, avg order
# value + avg delivery days recalculated accordingly).
# ---------------------------------------------------------------------------
$wsKPIs = $wb.Worksheets.Item(1)
$wsKPIs.Range("A2").Value = 6979.6
$wsKPIs.Range("E2").Value = 145.4083333333333
$wsKPIs.Range("F2").Value = 0.9204545454545454

# ---------------------------------------------------------------------------
# Sheet 2: "Ventes Mensuelles" - monthly sales totals + quantities revised.
# ---------------------------------------------------------------------------
$wsVentes = $wb.Worksheets.Item(2)
$wsVentes.Range("C2").Value = 681.2

$wsVentes.Range("C3").Value = 348.7
$wsVentes.Range("E3").Value = 22

$wsVentes.Range("C4").Value = 877.9
$wsVentes.Range("E4").Value = 57

$wsVentes.Range("C5").Value = 3025.7
$wsVentes.Range("E5").Value = 202

$wsVentes.Range("C6").Value = 994.9
$wsVentes.Range("E6").Value = 68

$wsVentes.Range("C7").Value = 1051.2
$wsVentes.Range("E7").Value = 76

# ---------------------------------------------------------------------------
# Sheet 3: "Par Catégorie" - Beverages now outranks Condiments, rows reorder.
# ---------------------------------------------------------------------------
$wsCat = $wb.Worksheets.Item(3)
$wsCat.Range("A2").Value = "Beverages"
$wsCat.Range("B2").Value = 4362.3
$wsCat.Range("C2").Value = 48
$wsCat.Range("D2").Value = 262

$wsCat.Range("A3").Value = "Condiments"
$wsCat.Range("B3").Value = 2617.3
$wsCat.Range("C3").Value = 28
$wsCat.Range("D3").Value = 208

# ---------------------------------------------------------------------------
# Sheet 4: "Top Produits" - Chai's numbers nudge slightly; Syrup now beats
# Cajun Seasoning for second place, rows reorder accordingly.
# ---------------------------------------------------------------------------
$wsProd = $wb.Worksheets.Item(4)
$wsProd.Range("B2").Value = 4362.3
$wsProd.Range("C2").Value = 262

$wsProd.Range("A3").Value = "Northwind Traders Syrup"
$wsProd.Range("B3").Value = 1332.5
$wsProd.Range("C3").Value = 145
$wsProd.Range("D3").Value = 28

$wsProd.Range("A4").Value = "Northwind Traders Cajun Seasoning"
$wsProd.Range("B4").Value = 1284.8
$wsProd.Range("C4").Value = 63
$wsProd.Range("D4").Value = 12

# ---------------------------------------------------------------------------
# Sheet 5: "Par Pays" - USA total revenue matches the new KPI figure.
# ---------------------------------------------------------------------------
$wsPays = $wb.Worksheets.Item(5)
$wsPays.Range("B2").Value = 6979.6

# ---------------------------------------------------------------------------
# Sheet 6: "Employés" - sales revised across the board; 2nd/3rd place swap
# (Mariya Sergienko now ahead of Anne Hellung-Larsen) and 7th/8th/9th place
# reshuffle (Andrew Cencini, Robert Zare, Laura Giussani).
# ---------------------------------------------------------------------------
$wsEmp = $wb.Worksheets.Item(6)
$wsEmp.Range("B2").Value = 1598.1

$wsEmp.Range("A3").Value = "Mariya Sergienko"
$wsEmp.Range("B3").Value = 1526.1
$wsEmp.Range("C3").Value = 8
$wsEmp.Range("D3").Value = 4

$wsEmp.Range("A4").Value = "Anne Hellung-Larsen"
$wsEmp.Range("B4").Value = 1332.2
$wsEmp.Range("C4").Value = 10
$wsEmp.Range("D4").Value = 5

$wsEmp.Range("B5").Value = 945.9

$wsEmp.Range("B6").Value = 701.2

$wsEmp.Range("A7").Value = "Andrew Cencini"
$wsEmp.Range("B7").Value = 419.8
$wsEmp.Range("C7").Value = 4
$wsEmp.Range("D7").Value = 3

$wsEmp.Range("A8").Value = "Robert Zare"
$wsEmp.Range("B8").Value = 254.7
$wsEmp.Range("C8").Value = 2
$wsEmp.Range("D8").Value = 1

$wsEmp.Range("A9").Value = "Laura Giussani"
$wsEmp.Range("B9").Value = 201.6
$wsEmp.Range("C9").Value = 2
$wsEmp.Range("D9").Value = 1
